$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NOTE row 11: enhancer.probes.Rd documentation note
$ws.Range("A11").Value = "NOTE"
$ws.Range("B11").Value = "enhancer.probes.Rd:10-12: Dropping empty section \source"

# New NOTE row 12: regions.Rd documentation note
$ws.Range("A12").Value = "NOTE"
$ws.Range("B12").Value = "regions.Rd:13-15: Dropping empty section \source"

# Match the "Comment" column formatting (Monaco font) used by the rest of column B
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B11:B12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection on the newly added row, matching the saved view state
$ws.Range("C12").Select() | Out-Null
